$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.778.09"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.798.01"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.15"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.50"
$ws.Range("E6").Value = "  +5.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.17"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.87"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.76"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("D15").Value = "3.236.65"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.966"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").Value = "2.789.51"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "51.754.29"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.36"
$ws.Range("E19").Value = "  +10.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  +4.68%  "
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.38"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.94"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.11"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.77"
$ws.Range("E29").Value = "  +13.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.44"
$ws.Range("E30").Value = "  +3.70%  "
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.71"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.10"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0913"
$ws.Range("E34").Value = "  +10.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0454"
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("E36").Value = "  +5.96%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.89"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.29"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.03"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +9.11%  "
$ws.Range("E47").Value = "  +8.85%  "
$ws.Range("D48").Value = "2.121.77"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.980"
$ws.Range("E49").Value = "  +7.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.50"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("E51").Value = "  +16.97%  "
